$d = $word.ActiveDocument

# --- Paragraph 1: "This is a Microsoft word document." ---
# Add trailing two spaces to the existing run, then append three new
# red-colored runs forming "(This is a change – Version for main branch)"
# split across three runs (mirrors the authored diff's run boundaries).
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs.Item(1)

$r1 = $p1.Range
$r1.MoveEnd(1, -1)
$r1.Collapse(0)
$r1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r1.Font.Color = 255

$r2 = $p1.Range
$r2.MoveEnd(1, -1)
$r2.Collapse(0)
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255

$r3 = $p1.Range
$r3.MoveEnd(1, -1)
$r3.Collapse(0)
$r3.InsertAfter(")")
$r3.Font.Color = 255

# --- Paragraph 3: the empty Menlo-styled paragraph becomes a bare <w:p/> ---
$p3 = $d.Paragraphs.Item(3)
$emptyParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$p3.Range.InsertXML($emptyParaXml)

Write-Output "edit complete"
